# Scheduled data refresh: update market-board price/profit snapshots
# across the per-job Profits sheets (commit: "chore: update Sheets via scheduled runner").
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 143.33333
$ws.Range("J4").Value = 80
$ws.Range("L4").Value = 80
$ws.Range("N4").Value = -308
# Row 9
$ws.Range("H9").Value = 190
$ws.Range("I9").Value = 180
$ws.Range("K9").Value = 180
$ws.Range("M9").Value = -11
# Row 11
$ws.Range("H11").Value = 1215.3334
$ws.Range("I11").Value = 1215.3334
$ws.Range("K11").Value = 1215.3334
$ws.Range("M11").Value = -1075.3334
# Row 33
$ws.Range("H33").Value = 117.29412
$ws.Range("I33").Value = 107.23077
$ws.Range("J33").Value = 150
$ws.Range("K33").Value = 107.23077
$ws.Range("L33").Value = 150
$ws.Range("M33").Value = 121.76923
$ws.Range("N33").Value = -608
# Row 107
$ws.Range("H107").Value = 221.16667
$ws.Range("I107").Value = 187
$ws.Range("K107").Value = 187
$ws.Range("M107").Value = 1733
# Row 113
$ws.Range("H113").Value = 3199.5
$ws.Range("J113").Value = 3199.5
$ws.Range("L113").Value = 3199.5
$ws.Range("N113").Value = -9707.5
# Row 116
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
# Row 137
$ws.Range("H137").Value = 2177.5
$ws.Range("I137").Value = 965.5
$ws.Range("J137").Value = 2480.5
$ws.Range("K137").Value = 2896.5
$ws.Range("L137").Value = 7441.5
$ws.Range("M137").Value = -346.5
$ws.Range("N137").Value = -12541.5
# Row 141
$ws.Range("H141").Value = 2327.1333
$ws.Range("I141").Value = 2327.1333
$ws.Range("K141").Value = 6981.3999
$ws.Range("M141").Value = -1801.3999

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 859.5
$ws.Range("I2").Value = 791.4
$ws.Range("J2").Value = 1200
$ws.Range("K2").Value = 791.4
$ws.Range("L2").Value = 1200
$ws.Range("M2").Value = -678.4
$ws.Range("N2").Value = -1426
# Row 45
$ws.Range("H45").Value = 2148.1667
$ws.Range("I45").Value = 1796.3334
$ws.Range("J45").Value = 2500
$ws.Range("K45").Value = 1796.3334
$ws.Range("L45").Value = 2500
$ws.Range("M45").Value = -1419.3334
$ws.Range("N45").Value = -3254
# Row 74
$ws.Range("H74").Value = 3414.9333
$ws.Range("I74").Value = 3473.1428
$ws.Range("K74").Value = 3473.1428
$ws.Range("M74").Value = -2599.1428
# Row 77
$ws.Range("H77").Value = 3414.9333
$ws.Range("I77").Value = 3473.1428
$ws.Range("K77").Value = 17365.714
$ws.Range("M77").Value = -12997.714
# Row 102
$ws.Range("H102").Value = 1436.7778
$ws.Range("I102").Value = 1191.375
$ws.Range("K102").Value = 1191.375
$ws.Range("M102").Value = 430.625
# Row 116
$ws.Range("H116").Value = 859.5
$ws.Range("I116").Value = 791.4
$ws.Range("J116").Value = 1200
$ws.Range("K116").Value = 791.4
$ws.Range("L116").Value = 1200
$ws.Range("M116").Value = 1502.6
$ws.Range("N116").Value = -5788
# Row 122
$ws.Range("H122").Value = 1448.5
$ws.Range("I122").Value = 1448.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4345.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1895.5
$ws.Range("N122").ClearContents()
# Row 132
$ws.Range("H132").Value = 1529.7727
$ws.Range("I132").Value = 1425.3334
$ws.Range("K132").Value = 4276.0002
$ws.Range("M132").Value = -1746.0002

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 859.5
$ws.Range("I3").Value = 791.4
$ws.Range("J3").Value = 1200
$ws.Range("K3").Value = 791.4
$ws.Range("L3").Value = 1200
$ws.Range("M3").Value = -677.4
$ws.Range("N3").Value = -1428
# Row 20
$ws.Range("H20").Value = 4008.6667
$ws.Range("I20").Value = 3265
$ws.Range("J20").Value = 5496
$ws.Range("K20").Value = 3265
$ws.Range("L20").Value = 5496
$ws.Range("M20").Value = -3018
$ws.Range("N20").Value = -5990
# Row 99
$ws.Range("H99").Value = 7600
$ws.Range("I99").Value = 7600
$ws.Range("K99").Value = 7600
$ws.Range("M99").Value = -6102
# Row 107
$ws.Range("H107").Value = 1136
$ws.Range("I107").Value = 1104
$ws.Range("K107").Value = 1104
$ws.Range("M107").Value = 816
# Row 134
$ws.Range("H134").Value = 3862
$ws.Range("I134").Value = 3862
$ws.Range("K134").Value = 11586
$ws.Range("M134").Value = -9051

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1814.772
$ws.Range("I31").Value = 1374.262
$ws.Range("J31").Value = 3048.2
$ws.Range("K31").Value = 1374.262
$ws.Range("L31").Value = 3048.2
$ws.Range("M31").Value = -1079.262
$ws.Range("N31").Value = -3638.2
# Row 34
$ws.Range("H34").Value = 1814.772
$ws.Range("I34").Value = 1374.262
$ws.Range("J34").Value = 3048.2
$ws.Range("K34").Value = 1374.262
$ws.Range("L34").Value = 3048.2
$ws.Range("M34").Value = -1172.262
$ws.Range("N34").Value = -3452.2
# Row 105
$ws.Range("H105").Value = 1396.3334
$ws.Range("I105").Value = 1219.5
$ws.Range("K105").Value = 1219.5
$ws.Range("M105").Value = 527.5
# Row 122
$ws.Range("H122").Value = 1017.8
$ws.Range("I122").Value = 897.25
$ws.Range("K122").Value = 2691.75
$ws.Range("M122").Value = -241.75
# Row 131
$ws.Range("H131").Value = 49999.5
$ws.Range("J131").Value = 49999.5
$ws.Range("L131").Value = 49999.5
$ws.Range("N131").Value = -60079.5

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 39
$ws.Range("H39").Value = 4955.4443
$ws.Range("J39").Value = 4955.4443
$ws.Range("L39").Value = 14866.3329
$ws.Range("N39").Value = -15454.3329
# Row 55
$ws.Range("H55").Value = 9202.474
$ws.Range("I55").Value = 14385
$ws.Range("J55").Value = 3444.111
$ws.Range("K55").Value = 43155
$ws.Range("L55").Value = 10332.333
$ws.Range("M55").Value = -42978
$ws.Range("N55").Value = -10686.333
# Row 118
$ws.Range("H118").Value = 2229.3333
$ws.Range("J118").Value = 6666
$ws.Range("L118").Value = 19998
$ws.Range("N118").Value = -22484
# Row 139
$ws.Range("H139").Value = 3378
$ws.Range("I139").Value = 854.8889
$ws.Range("K139").Value = 2564.6667
$ws.Range("M139").Value = 2575.3333

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2706.5625
$ws.Range("I80").Value = 2468.625
$ws.Range("J80").Value = 2944.5
$ws.Range("K80").Value = 2468.625
$ws.Range("L80").Value = 2944.5
$ws.Range("M80").Value = -1470.625
$ws.Range("N80").Value = -4940.5
# Row 83
$ws.Range("H83").Value = 2706.5625
$ws.Range("I83").Value = 2468.625
$ws.Range("J83").Value = 2944.5
$ws.Range("K83").Value = 12343.125
$ws.Range("L83").Value = 14722.5
$ws.Range("M83").Value = -7351.125
$ws.Range("N83").Value = -24706.5
# Row 97
$ws.Range("H97").Value = 635
$ws.Range("I97").Value = 759.8
$ws.Range("J97").Value = 11
$ws.Range("K97").Value = 759.8
$ws.Range("L97").Value = 11
$ws.Range("M97").Value = -263.8
$ws.Range("N97").Value = -1003
# Row 122
$ws.Range("H122").Value = 2745.5
$ws.Range("I122").Value = 1743.625
$ws.Range("J122").Value = 6753
$ws.Range("K122").Value = 5230.875
$ws.Range("L122").Value = 20259
$ws.Range("M122").Value = -2780.875
$ws.Range("N122").Value = -25159

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 123
$ws.Range("H123").Value = 78996
$ws.Range("J123").Value = 78996
$ws.Range("L123").Value = 78996
$ws.Range("N123").Value = -88796

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 59
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
# Row 62
$ws.Range("H62").Value = 5291.4546
$ws.Range("J62").Value = 5682.125
$ws.Range("L62").Value = 5682.125
$ws.Range("N62").Value = -6930.125
# Row 65
$ws.Range("H65").Value = 5291.4546
$ws.Range("J65").Value = 5682.125
$ws.Range("L65").Value = 28410.625
$ws.Range("N65").Value = -34650.625
# Row 132
$ws.Range("H132").Value = 843.7143
$ws.Range("I132").Value = 734.4167
$ws.Range("J132").Value = 1499.5
$ws.Range("K132").Value = 2203.2501
$ws.Range("L132").Value = 4498.5
$ws.Range("M132").Value = 326.7498999999998
$ws.Range("N132").Value = -9558.5
